$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value  = "Ki"
$ws.Range("A4").Value  = "ti"
$ws.Range("A5").Value  = "Firma2"
$ws.Range("A6").Value  = "Textfield-1"
$ws.Range("A7").Value  = "Representante del solicitante"
$ws.Range("A8").Value  = "h"
$ws.Range("A9").Value  = "NIFNIE"
$ws.Range("A10").Value = "Kp"
$ws.Range("A11").Value = "S"
$ws.Range("A13").Value = "Fecha fin actuación"
$ws.Range("A14").Value = "Fecha inicio actuación"
